$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Learn Python the hard way")

# Ex 40 - modules, class and objects
$ws.Range("B42").Value = 40
$ws.Range("C42").Value = "modules, class and objects"

# Ex 41 - OOP
$ws.Range("B43").Value = 41
$ws.Range("C43").Value = "OOP"

# Ex 42 - no notes
$ws.Range("B44").Value = 42

# Ex 43 - analysis and design
$ws.Range("B45").Value = 43
$ws.Range("C45").Value = "analysis and design"

# Ex 44 - inheritance and compostion
$ws.Range("B46").Value = 44
$ws.Range("C46").Value = "inheritance and compostion"

# Ex 45-51 - no notes yet
$ws.Range("B47").Value = 45
$ws.Range("B48").Value = 46
$ws.Range("B49").Value = 47
$ws.Range("B50").Value = 48
$ws.Range("B51").Value = 49
$ws.Range("B52").Value = 50
$ws.Range("B53").Value = 51

# Match the formatting applied to C42/C43 in the source workbook (same
# style as the surrounding "Ex" rows, e.g. C41).
$ws.Range("C41").Copy()
$ws.Range("C42:C43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the view where the author left it: scrolled down near the new
# rows with F52 selected.
$ws.Activate()
$ws.Range("F52").Select()
